$wb = $excel.ActiveWorkbook

# This script refreshes cached market-price / profit values in the
# FFXIV Leve profit-tracking sheets (Lamia server), mirroring a scheduled
# data-refresh run. Columns: H=currentAveragePrice, I=currentAveragePriceNQ,
# J=currentAveragePriceHQ, K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ,
# N=LeveProfitHQ. All values are static (no formulas in the workbook), so
# each changed cell is written explicitly.

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 5785.9287
$ws.Cells.Item(43, 9).Value = 4000.3333
$ws.Cells.Item(43, 10).Value = 6272.909
$ws.Cells.Item(43, 11).Value = 4000.3333
$ws.Cells.Item(43, 12).Value = 6272.909
$ws.Cells.Item(43, 13).Value = -3931.3333
$ws.Cells.Item(43, 14).Value = -6410.909
$ws.Cells.Item(64, 8).Value = 5749.5
$ws.Cells.Item(64, 10).Value = 4998
$ws.Cells.Item(64, 12).Value = 4998
$ws.Cells.Item(64, 14).Value = -5494
$ws.Cells.Item(67, 8).Value = 5749.5
$ws.Cells.Item(67, 10).Value = 4998
$ws.Cells.Item(67, 12).Value = 4998
$ws.Cells.Item(67, 14).Value = -6714
$ws.Cells.Item(132, 8).Value = 1561.2812
$ws.Cells.Item(132, 9).Value = 1387.3704
$ws.Cells.Item(132, 10).Value = 2500.4
$ws.Cells.Item(132, 11).Value = 4162.1112
$ws.Cells.Item(132, 12).Value = 7501.200000000001
$ws.Cells.Item(132, 13).Value = -1632.1112
$ws.Cells.Item(132, 14).Value = -12561.2
$ws.Cells.Item(137, 8).Value = 3538.6
$ws.Cells.Item(137, 9).Value = 2840.2
$ws.Cells.Item(137, 10).Value = 3771.4
$ws.Cells.Item(137, 11).Value = 8520.599999999999
$ws.Cells.Item(137, 12).Value = 11314.2
$ws.Cells.Item(137, 13).Value = -5970.599999999999
$ws.Cells.Item(137, 14).Value = -16414.2
$ws.Cells.Item(138, 8).Value = 3235.1636
$ws.Cells.Item(138, 10).Value = 3956.5625
$ws.Cells.Item(138, 12).Value = 11869.6875
$ws.Cells.Item(138, 14).Value = -22149.6875

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(24, 8).Value = 28000
$ws.Cells.Item(24, 10).Value = 28000
$ws.Cells.Item(24, 12).Value = 28000
$ws.Cells.Item(24, 14).Value = -28748
$ws.Cells.Item(100, 8).Value = 28000
$ws.Cells.Item(100, 10).Value = 28000
$ws.Cells.Item(100, 12).Value = 28000
$ws.Cells.Item(100, 14).Value = -30164
$ws.Cells.Item(101, 8).Value = 29998.75
$ws.Cells.Item(101, 10).Value = 29998.75
$ws.Cells.Item(101, 12).Value = 29998.75
$ws.Cells.Item(101, 14).Value = -36488.75
$ws.Cells.Item(102, 8).Value = 1584.0555
$ws.Cells.Item(102, 10).Value = 1100
$ws.Cells.Item(102, 12).Value = 1100
$ws.Cells.Item(102, 14).Value = -4344
$ws.Cells.Item(114, 8).Value = 69999
$ws.Cells.Item(114, 10).Value = 69999
$ws.Cells.Item(114, 12).Value = 69999
$ws.Cells.Item(114, 14).Value = -78677
$ws.Cells.Item(115, 8).Value = 9995
$ws.Cells.Item(115, 10).Value = 9995
$ws.Cells.Item(115, 12).Value = 9995
$ws.Cells.Item(115, 14).Value = -13129

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(6, 8).Value = 9999.625
$ws.Cells.Item(6, 10).Value = 9999.625
$ws.Cells.Item(6, 12).Value = 9999.625
$ws.Cells.Item(6, 14).Value = -10225.625
$ws.Cells.Item(86, 8).Value = 2812.6316
$ws.Cells.Item(86, 9).Value = 1995.6154
$ws.Cells.Item(86, 10).Value = 4582.8335
$ws.Cells.Item(86, 11).Value = 1995.6154
$ws.Cells.Item(86, 12).Value = 4582.8335
$ws.Cells.Item(86, 13).Value = -872.6153999999999
$ws.Cells.Item(86, 14).Value = -6828.8335
$ws.Cells.Item(89, 8).Value = 2812.6316
$ws.Cells.Item(89, 9).Value = 1995.6154
$ws.Cells.Item(89, 10).Value = 4582.8335
$ws.Cells.Item(89, 11).Value = 9978.076999999999
$ws.Cells.Item(89, 12).Value = 22914.1675
$ws.Cells.Item(89, 13).Value = -4362.076999999999
$ws.Cells.Item(89, 14).Value = -34146.1675
$ws.Cells.Item(105, 8).Value = 13250.3
$ws.Cells.Item(105, 9).Value = 6917.3335
$ws.Cells.Item(105, 10).Value = 22749.75
$ws.Cells.Item(105, 11).Value = 6917.3335
$ws.Cells.Item(105, 12).Value = 22749.75
$ws.Cells.Item(105, 13).Value = -5170.3335
$ws.Cells.Item(105, 14).Value = -26243.75
$ws.Cells.Item(114, 8).Value = 43200
$ws.Cells.Item(114, 10).Value = 43200
$ws.Cells.Item(114, 12).Value = 43200
$ws.Cells.Item(114, 14).Value = -51878
$ws.Cells.Item(117, 8).Value = 73684
$ws.Cells.Item(117, 10).Value = 73684
$ws.Cells.Item(117, 12).Value = 73684
$ws.Cells.Item(117, 14).Value = -82862
$ws.Cells.Item(119, 8).Value = 68421
$ws.Cells.Item(119, 10).Value = 68421
$ws.Cells.Item(119, 12).Value = 68421
$ws.Cells.Item(119, 14).Value = -78097
$ws.Cells.Item(120, 8).Value = 19999.5
$ws.Cells.Item(120, 10).Value = 19999.5
$ws.Cells.Item(120, 12).Value = 19999.5
$ws.Cells.Item(120, 14).Value = -29675.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 3502.4285
$ws.Cells.Item(16, 9).Value = 1999.25
$ws.Cells.Item(16, 10).Value = 5506.6665
$ws.Cells.Item(16, 11).Value = 1999.25
$ws.Cells.Item(16, 12).Value = 5506.6665
$ws.Cells.Item(16, 13).Value = -1712.25
$ws.Cells.Item(16, 14).Value = -6080.6665
$ws.Cells.Item(99, 8).Value = 2833.3333
$ws.Cells.Item(99, 9).Value = 2750
$ws.Cells.Item(99, 11).Value = 2750
$ws.Cells.Item(99, 13).Value = -1252
$ws.Cells.Item(105, 8).Value = 10901.333
$ws.Cells.Item(105, 9).Value = 6849.5
$ws.Cells.Item(105, 10).Value = 19005
$ws.Cells.Item(105, 11).Value = 6849.5
$ws.Cells.Item(105, 12).Value = 19005
$ws.Cells.Item(105, 13).Value = -5102.5
$ws.Cells.Item(105, 14).Value = -22499
$ws.Cells.Item(113, 8).Value = 3502.4285
$ws.Cells.Item(113, 9).Value = 1999.25
$ws.Cells.Item(113, 10).Value = 5506.6665
$ws.Cells.Item(113, 11).Value = 1999.25
$ws.Cells.Item(113, 12).Value = 5506.6665
$ws.Cells.Item(113, 13).Value = 170.75
$ws.Cells.Item(113, 14).Value = -9846.666499999999
$ws.Cells.Item(114, 8).Value = 46250
$ws.Cells.Item(114, 10).Value = 46250
$ws.Cells.Item(114, 12).Value = 46250
$ws.Cells.Item(114, 14).Value = -54928
$ws.Cells.Item(126, 8).Value = 2833.3333
$ws.Cells.Item(126, 9).Value = 2750
$ws.Cells.Item(126, 11).Value = 8250
$ws.Cells.Item(126, 13).Value = -5780
$ws.Cells.Item(132, 8).Value = 2821.0557
$ws.Cells.Item(132, 9).Value = 3142.7144
$ws.Cells.Item(132, 10).Value = 2370.7334
$ws.Cells.Item(132, 11).Value = 9428.143199999999
$ws.Cells.Item(132, 12).Value = 7112.2002
$ws.Cells.Item(132, 13).Value = -6898.143199999999
$ws.Cells.Item(132, 14).Value = -12172.2002
$ws.Cells.Item(134, 8).Value = 1771.7317
$ws.Cells.Item(134, 9).Value = 1143.0646
$ws.Cells.Item(134, 11).Value = 3429.1938
$ws.Cells.Item(134, 13).Value = -894.1938

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 7645847
$ws.Cells.Item(4, 9).Value = 4666765
$ws.Cells.Item(4, 11).Value = 14000295
$ws.Cells.Item(4, 13).Value = -14000183
$ws.Cells.Item(60, 8).Value = 912651
$ws.Cells.Item(60, 9).Value = 379
$ws.Cells.Item(60, 11).Value = 1137
$ws.Cells.Item(60, 13).Value = -886
$ws.Cells.Item(68, 8).Value = 4418.75
$ws.Cells.Item(68, 10).Value = 4478.5713
$ws.Cells.Item(68, 12).Value = 13435.7139
$ws.Cells.Item(68, 14).Value = -15057.7139
$ws.Cells.Item(71, 8).Value = 4418.75
$ws.Cells.Item(71, 10).Value = 4478.5713
$ws.Cells.Item(71, 12).Value = 40307.14169999999
$ws.Cells.Item(71, 14).Value = -48419.14169999999
$ws.Cells.Item(81, 8).Value = 2549.75
$ws.Cells.Item(81, 9).Value = 2399.6667
$ws.Cells.Item(81, 10).Value = 3000
$ws.Cells.Item(81, 11).Value = 7199.000100000001
$ws.Cells.Item(81, 12).Value = 9000
$ws.Cells.Item(81, 13).Value = -6076.000100000001
$ws.Cells.Item(81, 14).Value = -11246
$ws.Cells.Item(84, 8).Value = 2549.75
$ws.Cells.Item(84, 9).Value = 2399.6667
$ws.Cells.Item(84, 10).Value = 3000
$ws.Cells.Item(84, 11).Value = 21597.0003
$ws.Cells.Item(84, 12).Value = 27000
$ws.Cells.Item(84, 13).Value = -15981.0003
$ws.Cells.Item(84, 14).Value = -38232
$ws.Cells.Item(103, 8).Value = 1375.9
$ws.Cells.Item(103, 9).Value = 835
$ws.Cells.Item(103, 10).Value = 2187.25
$ws.Cells.Item(103, 11).Value = 2505
$ws.Cells.Item(103, 12).Value = 6561.75
$ws.Cells.Item(103, 13).Value = -1626
$ws.Cells.Item(103, 14).Value = -8319.75
$ws.Cells.Item(118, 8).Value = 2820
$ws.Cells.Item(118, 9).Value = 980
$ws.Cells.Item(118, 11).Value = 2940
$ws.Cells.Item(118, 13).Value = -1697
$ws.Cells.Item(132, 8).Value = 4862.385
$ws.Cells.Item(132, 10).Value = 4521.1
$ws.Cells.Item(132, 12).Value = 40689.9
$ws.Cells.Item(132, 14).Value = -45749.9

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(7, 8).Value = 25249.5
$ws.Cells.Item(7, 10).Value = 27000
$ws.Cells.Item(7, 12).Value = 27000
$ws.Cells.Item(7, 14).Value = -27224
$ws.Cells.Item(8, 8).Value = 25249.5
$ws.Cells.Item(8, 10).Value = 27000
$ws.Cells.Item(8, 12).Value = 27000
$ws.Cells.Item(8, 14).Value = -27278
$ws.Cells.Item(70, 8).Value = 16026.35
$ws.Cells.Item(70, 9).Value = 6035.467
$ws.Cells.Item(70, 10).Value = 45999
$ws.Cells.Item(70, 11).Value = 6035.467
$ws.Cells.Item(70, 12).Value = 45999
$ws.Cells.Item(70, 13).Value = -5765.467
$ws.Cells.Item(70, 14).Value = -46539
$ws.Cells.Item(73, 8).Value = 16026.35
$ws.Cells.Item(73, 9).Value = 6035.467
$ws.Cells.Item(73, 10).Value = 45999
$ws.Cells.Item(73, 11).Value = 6035.467
$ws.Cells.Item(73, 12).Value = 45999
$ws.Cells.Item(73, 13).Value = -5099.467
$ws.Cells.Item(73, 14).Value = -47871
$ws.Cells.Item(97, 8).Value = 937.4737
$ws.Cells.Item(97, 9).Value = 788.6667
$ws.Cells.Item(97, 10).Value = 1071.4
$ws.Cells.Item(97, 11).Value = 788.6667
$ws.Cells.Item(97, 12).Value = 1071.4
$ws.Cells.Item(97, 13).Value = -292.6667
$ws.Cells.Item(97, 14).Value = -2063.4
$ws.Cells.Item(106, 8).Value = 69264
$ws.Cells.Item(106, 10).Value = 69264
$ws.Cells.Item(106, 12).Value = 69264
$ws.Cells.Item(106, 14).Value = -71788
$ws.Cells.Item(107, 8).Value = 977.16
$ws.Cells.Item(107, 9).Value = 405.76923
$ws.Cells.Item(107, 10).Value = 1596.1666
$ws.Cells.Item(107, 11).Value = 405.76923
$ws.Cells.Item(107, 12).Value = 1596.1666
$ws.Cells.Item(107, 13).Value = 1514.23077
$ws.Cells.Item(107, 14).Value = -5436.1666
$ws.Cells.Item(126, 8).Value = 7088.6
$ws.Cells.Item(126, 9).Value = 4489.5
$ws.Cells.Item(126, 10).Value = 8821.333000000001
$ws.Cells.Item(126, 11).Value = 13468.5
$ws.Cells.Item(126, 12).Value = 26463.999
$ws.Cells.Item(126, 13).Value = -10998.5
$ws.Cells.Item(126, 14).Value = -31403.999
$ws.Cells.Item(136, 8).Value = 43661.355
$ws.Cells.Item(136, 10).Value = 43661.355
$ws.Cells.Item(136, 12).Value = 130984.065
$ws.Cells.Item(136, 14).Value = -136084.065
